# Generate Report for Handback
# The bf579bfb-c5ae-4046-a4b7-e4ceda83f2ac file has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales, so update the
# Overview sheet plus each locale's detail sheet with the new status and
# handback timestamps, and clear the stale error detail.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the bf579bfb-c5ae-4046-a4b7-e4ceda83f2ac file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the bf579bfb-c5ae-4046-a4b7-e4ceda83f2ac file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-11-14 05:51:55"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row 3 is the bf579bfb-c5ae-4046-a4b7-e4ceda83f2ac file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-11-14 05:52:14"
$wsDeDe.Range("P3").Value = ""
